$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C and D contain plain text in the source workbook (coin names,
# links and price strings such as "1.00" or "0.539"). Assigning those through
# .Value would let Excel auto-detect them as numbers and silently reformat
# them (e.g. "1.00" -> 1). Temporarily force the cell to Text format while
# writing the new value, then restore the original style so no formatting
# is left behind.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "31.597.73"
$ws.Range("E2").Value = "  +5.91%  "
Set-TextValue $ws.Range("D3") "1.719.24"
$ws.Range("E3").Value = "  +4.96%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue $ws.Range("D5") "223.81"
$ws.Range("E5").Value = "  +4.01%  "
Set-TextValue $ws.Range("D6") "0.539"
$ws.Range("E6").Value = "  +3.80%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.13%  "
Set-TextValue $ws.Range("D8") "30.04"
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("E9").Value = "  +3.39%  "
Set-TextValue $ws.Range("D10") "0.0654"
$ws.Range("E10").Value = "  +7.36%  "
Set-TextValue $ws.Range("D11") "0.0910"
$ws.Range("E11").Value = "  +1.34%  "
Set-TextValue $ws.Range("D12") "1.962.01"
$ws.Range("E12").Value = "  +4.76%  "
Set-TextValue $ws.Range("D13") "1.708.46"
$ws.Range("E13").Value = "  +4.20%  "
$ws.Range("E14").Value = "  +4.27%  "
Set-TextValue $ws.Range("D15") "10.17"
$ws.Range("E15").Value = "  +6.91%  "
$ws.Range("E16").Value = "  +7.69%  "
Set-TextValue $ws.Range("D17") "31.604.99"
$ws.Range("E17").Value = "  +5.89%  "
$ws.Range("E18").Value = "  +4.92%  "
Set-TextValue $ws.Range("D19") "251.66"
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("E20").Value = "  +3.12%  "
Set-TextValue $ws.Range("D21") "0.999"
$ws.Range("E21").Value = "  -0.11%  "
Set-TextValue $ws.Range("D22") "10.16"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("E23").Value = "  +3.04%  "
Set-TextValue $ws.Range("D24") "2.17"
$ws.Range("E24").Value = "  -0.16%  "
Set-TextValue $ws.Range("D25") "159.17"
$ws.Range("E25").Value = "  +1.25%  "
Set-TextValue $ws.Range("D26") "16.13"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("E27").Value = "  +3.58%  "
Set-TextValue $ws.Range("D28") "6.83"
$ws.Range("E28").Value = "  +3.07%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.06%  "
Set-TextValue $ws.Range("D30") "3.87"
$ws.Range("E30").Value = "  +14.37%  "
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +4.91%  "
Set-TextValue $ws.Range("D33") "3.41"
$ws.Range("E33").Value = "  +6.91%  "
Set-TextValue $ws.Range("D34") "1.529.65"
$ws.Range("E34").Value = "  +7.65%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +2.61%  "
Set-TextValue $ws.Range("D37") "83.16"
$ws.Range("E37").Value = "  +8.68%  "
Set-TextValue $ws.Range("D38") "0.613"
$ws.Range("E38").Value = "  +8.48%  "
Set-TextValue $ws.Range("D39") "0.0182"
$ws.Range("E39").Value = "  +5.10%  "
Set-TextValue $ws.Range("D40") "2.74"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +0.67%  "
Set-TextValue $ws.Range("B42") "ARBITRUM"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D42") "0.856"
$ws.Range("E42").Value = "  +2.90%  "
Set-TextValue $ws.Range("B43") "RenderToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D43") "2.04"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("E44").Value = "  +0.90%  "
Set-TextValue $ws.Range("D45") "1.03"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("E46").Value = "  -0.11%  "
Set-TextValue $ws.Range("D47") "52.70"
$ws.Range("E47").Value = "  +5.78%  "
$ws.Range("E48").Value = "  +5.19%  "
Set-TextValue $ws.Range("D49") "1.848.75"
$ws.Range("E50").Value = "  +8.83%  "
Set-TextValue $ws.Range("D51") "93.66"
$ws.Range("E51").Value = "  +0.45%  "
